# Automatische test-sync: 2025-06-22 22:20:50
#
# Adds the new "Productinformatie" mail-log entry (row 61) to the "Logs"
# sheet, extends the dimension / conditional-formatting ranges to cover it,
# and re-syncs the "Dashboard" summary table (category counts + sort order)
# to reflect the extra Productinformatie row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet - append the new row at the bottom (row 61)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$answer = "Beste klant,`n" + `
    "Dank u voor uw interesse in onze producten A en B. Product A heeft een grotere capaciteit en meer geavanceerde functies dan product B. Product B daarentegen is compacter en eenvoudiger in gebruik. Afhankelijk van uw behoeften en budget kunt u kiezen welk product het beste bij u past.`n" + `
    "Mocht u nog specifieke vragen hebben over de functies van beide producten, dan helpen wij u graag verder.`n" + `
    "Met vriendelijke groet,`n" + `
    "[Naam van het bedrijf] E-mailassistent"

$logs.Range("A61").Value = "Productinformatie"
$logs.Range("B61").Value = "mailmind.test@zohomail.eu"
$logs.Range("C61").Value = "Wat is het verschil tussen product A en product B?"
$logs.Range("D61").Value = "Productinformatie"
$logs.Range("E61").Value = $answer
$logs.Range("F61").Value = "2025-06-22 22:19:53"
$logs.Range("G61").Value = "Ja"

# Entering the multi-line answer makes the engine flag the row with an
# explicit (wrapped-text) height; the cell isn't word-wrapped though, so
# re-fitting the row collapses it back to the sheet's default row height,
# matching every other (un-styled) row in the log.
$logs.Rows.Item(61).EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting ranges from row 60 to row 61
#    (they don't grow automatically when a new row is populated).
# ---------------------------------------------------------------------
$catRules = $logs.Range("D2:D60").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D61"))
}

$answeredRules = $logs.Range("G2:G60").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G61"))
}

# ---------------------------------------------------------------------
# 3. "Dashboard" sheet - resync the category/count summary table.
#    Productinformatie's count goes 5 -> 6, which ties it with
#    "Sollicitatie / Vacature" (6) and reorders rows 4-6.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Productinformatie"
$dash.Range("B4").Value = 6

$dash.Range("A5").Value = "Sollicitatie / Vacature"
$dash.Range("B5").Value = 6

$dash.Range("A6").Value = "Offerte / Prijsaanvraag"
$dash.Range("B6").Value = 5
